# Fixed naive component forecaster bug - Presentation state 11.02.
# Recomputes the QoQ naive-forecast error diagonal (AVERAGE_1_9_qoq) for horizons
# that the buggy forecaster previously left at 0 / blank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24
$ws.Range("K24").Value = -4.27960679849052

# Row 25
$ws.Range("J25").Value = -3.788996755212807
$ws.Range("K25").Value = 1.143562595650536

# Row 26
$ws.Range("I26").Value = -3.865333711727925
$ws.Range("J26").Value = 1.067225639135418
$ws.Range("K26").Value = 2.06956845314902

# Row 27
$ws.Range("H27").Value = -4.656609552145667
$ws.Range("I27").Value = 0.2759497987176758
$ws.Range("J27").Value = 1.278292612731277
$ws.Range("K27").Value = -2.341393141007416

# Row 28
$ws.Range("G28").Value = -4.465333711727925
$ws.Range("H28").Value = 0.4672256391354182
$ws.Range("I28").Value = 1.46956845314902
$ws.Range("J28").Value = -2.150117300589673
$ws.Range("K28").Value = -1.831122027629067

# Row 29
$ws.Range("F29").Value = -4.177479830936775
$ws.Range("G29").Value = 0.7550795199265679
$ws.Range("H29").Value = 1.75742233394017
$ws.Range("I29").Value = -1.862263419798524
$ws.Range("J29").Value = -1.543268146837917
$ws.Range("K29").Value = 1.037611590612414

# Row 30
$ws.Range("E30").Value = -5.069663971681845
$ws.Range("F30").Value = -0.1371046208185024
$ws.Range("G30").Value = 0.8652381931950993
$ws.Range("H30").Value = -2.754447560543594
$ws.Range("I30").Value = -2.435452287582987
$ws.Range("J30").Value = 0.1454274498673439
$ws.Range("K30").Value = -2.17953185948812

# Row 31
$ws.Range("D31").Value = -3.946978334557187
$ws.Range("E31").Value = 0.9855810163061562
$ws.Range("F31").Value = 1.987923830319758
$ws.Range("G31").Value = -1.631761923418935
$ws.Range("H31").Value = -1.312766650458329
$ws.Range("I31").Value = 1.268113086992003
$ws.Range("J31").Value = -1.056846222363461
$ws.Range("K31").Value = -0.001587209089905173

# Row 32
$ws.Range("C32").Value = -6.165333711727925
$ws.Range("D32").Value = -1.232774360864582
$ws.Range("E32").Value = -0.2304315468509801
$ws.Range("F32").Value = -3.850117300589673
$ws.Range("G32").Value = -3.531122027629067
$ws.Range("H32").Value = -0.9502422901787355
$ws.Range("I32").Value = -3.275201599534199
$ws.Range("J32").Value = -2.219942586260643
$ws.Range("K32").Value = -2.001652513402433

# Row 33
$ws.Range("B33").Value = -11.89637784493493
$ws.Range("C33").Value = -6.963818494071589
$ws.Range("D33").Value = -5.961475680057987
$ws.Range("E33").Value = -9.58116143379668
$ws.Range("F33").Value = -9.262166160836074
$ws.Range("G33").Value = -6.681286423385743
$ws.Range("H33").Value = -9.006245732741206
$ws.Range("I33").Value = -7.95098671946765
$ws.Range("J33").Value = -7.73269664660944
$ws.Range("K33").Value = -6.641963848986791

# Row 34
$ws.Range("B34").Value = 4.932559350863343
$ws.Range("C34").Value = 5.934902164876945
$ws.Range("D34").Value = 2.315216411138252
$ws.Range("E34").Value = 2.634211684098858
$ws.Range("F34").Value = 5.215091421549189
$ws.Range("G34").Value = 2.890132112193726
$ws.Range("H34").Value = 3.945391125467282
$ws.Range("I34").Value = 4.163681198325492
$ws.Range("J34").Value = 5.254413995948141
$ws.Range("K34").Value = 5.671137314123342

# Row 35
$ws.Range("B35").Value = 1.002342814013602
$ws.Range("C35").Value = -2.617342939725091
$ws.Range("D35").Value = -2.298347666764485
$ws.Range("E35").Value = 0.2825320706858463
$ws.Range("F35").Value = -2.042427238669617
$ws.Range("G35").Value = -0.9871682253960614
$ws.Range("H35").Value = -0.7688781525378516
$ws.Range("I35").Value = 0.3218546450847981
$ws.Range("J35").Value = 0.7385779632599989
$ws.Range("K35").Value = -0.8462048329513294

# Row 36
$ws.Range("B36").Value = -3.619685753738693
$ws.Range("C36").Value = -3.300690480778087
$ws.Range("D36").Value = -0.7198107433277554
$ws.Range("E36").Value = -3.044770052683219
$ws.Range("F36").Value = -1.989511039409663
$ws.Range("G36").Value = -1.771220966551453
$ws.Range("H36").Value = -0.6804881689288036
$ws.Range("I36").Value = -0.2637648507536028
$ws.Range("J36").Value = -1.848547646964931
$ws.Range("K36").Value = -0.9013296832945341

# Row 37
$ws.Range("B37").Value = 0.3189952729606063
$ws.Range("C37").Value = 2.899875010410938
$ws.Range("D37").Value = 0.5749157010554742
$ws.Range("E37").Value = 1.63017471432903
$ws.Range("F37").Value = 1.84846478718724
$ws.Range("G37").Value = 2.939197584809889
$ws.Range("H37").Value = 3.35592090298509
$ws.Range("I37").Value = 1.771138106773762
$ws.Range("J37").Value = 2.718356070444159
$ws.Range("K37").Value = 2.63866197424673

# Row 38
$ws.Range("B38").Value = 2.580879737450331
$ws.Range("C38").Value = 0.2559204280948679
$ws.Range("D38").Value = 1.311179441368424
$ws.Range("E38").Value = 1.529469514226633
$ws.Range("F38").Value = 2.620202311849283
$ws.Range("G38").Value = 3.036925630024484
$ws.Range("H38").Value = 1.452142833813156
$ws.Range("I38").Value = 2.399360797483553
$ws.Range("J38").Value = 2.319666701286124
$ws.Range("K38").Value = 2.031122027629067

# Row 39
$ws.Range("B39").Value = -2.324959309355464
$ws.Range("C39").Value = -1.269700296081908
$ws.Range("D39").Value = -1.051410223223698
$ws.Range("E39").Value = 0.03932257439895182
$ws.Range("F39").Value = 0.4560458925741526
$ws.Range("G39").Value = -1.128736903637176
$ws.Range("H39").Value = -0.1815189399667787
$ws.Range("I39").Value = -0.2612130361642073
$ws.Range("J39").Value = -0.5497577098212645
$ws.Range("K39").Value = -0.5592990347223861

# Row 40
$ws.Range("B40").Value = 1.055259013273556
$ws.Range("C40").Value = 1.273549086131766
$ws.Range("D40").Value = 2.364281883754416
$ws.Range("E40").Value = 2.781005201929616
$ws.Range("F40").Value = 1.196222405718288
$ws.Range("G40").Value = 2.143440369388685
$ws.Range("H40").Value = 2.063746273191256
$ws.Range("I40").Value = 1.775201599534199
$ws.Range("J40").Value = 1.765660274633077
$ws.Range("K40").Value = 0.7463330244572151

# Row 41
$ws.Range("B41").Value = 0.2182900728582098
$ws.Range("C41").Value = 1.30902287048086
$ws.Range("D41").Value = 1.72574618865606
$ws.Range("E41").Value = 0.140963392444732
$ws.Range("F41").Value = 1.088181356115129
$ws.Range("G41").Value = 1.0084872599177
$ws.Range("H41").Value = 0.7199425862606432
$ws.Range("I41").Value = 0.7104012613595216
$ws.Range("J41").Value = -0.3089259888163407
$ws.Range("K41").Value = -0.2871171714278518

# Row 42
$ws.Range("B42").Value = 1.09073279762265
$ws.Range("C42").Value = 1.507456115797851
$ws.Range("D42").Value = -0.07732668041347779
$ws.Range("E42").Value = 0.8698912832569192
$ws.Range("F42").Value = 0.7901971870594906
$ws.Range("G42").Value = 0.5016525134024334
$ws.Range("H42").Value = 0.4921111885013119
$ws.Range("I42").Value = -0.5272160616745505
$ws.Range("J42").Value = -0.5054072442860615
$ws.Range("K42").Value = 1.033562986125958

# Row 43
$ws.Range("B43").Value = 0.4167233181752008
$ws.Range("C43").Value = -1.168059478036128
$ws.Range("D43").Value = -0.2208415143657305
$ws.Range("E43").Value = -0.3005356105631591
$ws.Range("F43").Value = -0.5890802842202163
$ws.Range("G43").Value = -0.5986216091213379
$ws.Range("H43").Value = -1.6179488592972
$ws.Range("I43").Value = -1.596140041908711
$ws.Range("J43").Value = -0.05716981149669209
$ws.Range("K43").Value = -1.009170151698186

# Row 44
$ws.Range("B44").Value = -1.584782796211328
$ws.Range("C44").Value = -0.6375648325409313
$ws.Range("D44").Value = -0.7172589287383599
$ws.Range("E44").Value = -1.005803602395417
$ws.Range("F44").Value = -1.015344927296539
$ws.Range("G44").Value = -2.034672177472401
$ws.Range("H44").Value = -2.012863360083912
$ws.Range("I44").Value = -0.4738931296718929
$ws.Range("J44").Value = -1.425893469873386

# Row 45
$ws.Range("B45").Value = 0.947217963670397
$ws.Range("C45").Value = 0.8675238674729684
$ws.Range("D45").Value = 0.5789791938159112
$ws.Range("E45").Value = 0.5694378689147896
$ws.Range("F45").Value = -0.4498893812610727
$ws.Range("G45").Value = -0.4280805638725838
$ws.Range("H45").Value = 1.110889666539435
$ws.Range("I45").Value = 0.1588893263379418

# Row 46
$ws.Range("B46").Value = -0.07969409619742857
$ws.Range("C46").Value = -0.3682387698544858
$ws.Range("D46").Value = -0.3777800947556074
$ws.Range("E46").Value = -1.39710734493147
$ws.Range("F46").Value = -1.375298527542981
$ws.Range("G46").Value = 0.1636717028690384
$ws.Range("H46").Value = -0.7883286373324552

# Row 47
$ws.Range("B47").Value = -0.2885446736570572
$ws.Range("C47").Value = -0.2980859985581788
$ws.Range("D47").Value = -1.317413248734041
$ws.Range("E47").Value = -1.295604431345552
$ws.Range("F47").Value = 0.243365799066467
$ws.Range("G47").Value = -0.7086345411350266

# Row 48
$ws.Range("B48").Value = -0.009541324901121584
$ws.Range("C48").Value = -1.028868575076984
$ws.Range("D48").Value = -1.007059757688495
$ws.Range("E48").Value = 0.5319104727235242
$ws.Range("F48").Value = -0.4200898674779694

# Row 49
$ws.Range("B49").Value = -1.019327250175862
$ws.Range("C49").Value = -0.9975184327873734
$ws.Range("D49").Value = 0.5414517976246458
$ws.Range("E49").Value = -0.4105485425768478

# Row 50
$ws.Range("B50").Value = 0.02180881738848894
$ws.Range("C50").Value = 1.560779047800508
$ws.Range("D50").Value = 0.6087787075990145

# Row 51
$ws.Range("B51").Value = 1.538970230412019
$ws.Range("C51").Value = 0.5869698902105256

# Row 52
$ws.Range("B52").Value = -0.9520003402014936
